$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 6644
$ws.Range("F6").Value = 671
$ws.Range("F7").Value = 217
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 1055
$ws.Range("F10").Value = 827
$ws.Range("F11").Value = 1002
$ws.Range("F12").Value = 1298
$ws.Range("F16").Value = 526
$ws.Range("F21").Value = 711
$ws.Range("F24").Value = 92
$ws.Range("F25").Value = 1114
$ws.Range("F26").Value = 232
$ws.Range("F27").Value = 2330
$ws.Range("F29").Value = 1173
$ws.Range("F32").Value = 3780
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 5
$ws.Range("F18").Value = 326
$ws.Range("F19").Value = 4115
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1236
$ws.Range("F5").Value = 1616
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1236
$ws.Range("F4").Value = 1616
$ws.Range("F9").Value = 6644
$ws.Range("F14").Value = 671
$ws.Range("F15").Value = 671
$ws.Range("F16").Value = 217
$ws.Range("F17").Value = 1055
$ws.Range("F18").Value = 827
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 1298
$ws.Range("F28").Value = 526
$ws.Range("F35").Value = 711
$ws.Range("F38").Value = 92
$ws.Range("F41").Value = 1114
$ws.Range("F42").Value = 232
$ws.Range("F43").Value = 2330
$ws.Range("F47").Value = 1173
$ws.Range("F49").Value = 3780
